$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.935.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.768.03"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4549"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.12%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3530"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.01"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07385"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.096"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.007"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.191"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.767.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.62"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.91%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06441"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.773"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.965.36"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.21"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.974.42"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.154"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.077"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09222"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.610"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.657"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.86"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02283"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06120"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2092"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.950"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6253"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.380"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.823"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.22"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.732"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5848"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.70"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.935"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06825"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.75%  "
